# Updated symbol list — refresh coin prices / volume labels and fix a
# three-row ordering mix-up (BKEXToken / CEJI / KickToken) on the "Sheet1"
# coin-ranking table.
#
# All of the "Price" column (D) entries are stored as TEXT (not numbers) in
# this workbook, e.g. "244.74" rather than 244.74 — so every numeric-looking
# write below is entered with a leading apostrophe, exactly like typing
# '244.61 into Excel, to keep the cell's value type as Text instead of
# letting Excel auto-convert it to a Number (which would also silently
# drop meaningful trailing zeros, e.g. "0.1050" -> 105).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Range("D2").Value = "'244.61"
# Row 3 - OKB
$ws.Range("D3").Value = "'23.93"
# Row 4 - HuobiToken
$ws.Range("D4").Value = "'5.205"
# Row 6 - KuCoinToken
$ws.Range("D6").Value = "'6.485"
# Row 8 - MXToken
$ws.Range("D8").Value = "'0.8138"
# Row 9 - FTXToken
$ws.Range("D9").Value = "'0.8692"
# Row 10 - WazirX
$ws.Range("D10").Value = "'0.1369"
# Row 11 - MandalaExchangeToken
$ws.Range("D11").Value = "'0.06935"
# Row 13 - BitrueCoin
$ws.Range("D13").Value = "'0.02926"
# Row 14 - BitMartToken
$ws.Range("D14").Value = "'0.09320"
# Row 15 - MCDex
$ws.Range("D15").Value = "'3.849"
# Row 16 - BitForexToken
$ws.Range("D16").Value = "'0.001545"
# Row 17 - CoinExToken
$ws.Range("D17").Value = "'0.04717"
# Row 18 - One
$ws.Range("D18").Value = "'0.0006012"
$ws.Range("E18").Value = "17OneONEWorstin24h"
# Row 19 - TigerCash
$ws.Range("D19").Value = "'0.006211"
# Row 20 - BitKan
$ws.Range("D20").Value = "'0.001242"
# Row 21 - HotbitToken
$ws.Range("D21").Value = "'0.004106"
# Row 22 - NitroEx
$ws.Range("D22").Value = "'0.00007005"
# Row 23 - LEO
$ws.Range("D23").Value = "'3.550"
# Row 25 - BitpandaEcosystemToken
$ws.Range("D25").Value = "'0.3192"
# Row 27 - UpBots
$ws.Range("D27").Value = "'0.0002329"
# Row 40 - IDEX
$ws.Range("D40").Value = "'0.03711"

# Rows 41-43 were re-ordered: what used to be KickToken/BKEXToken/CEJI (in
# that original relative order across rows 41/42/43 as BKEXToken, CEJI,
# KickToken) is now listed as KickToken, BKEXToken, CEJI, each with a
# refreshed price and rank-label.
# Row 41 - now KickToken (was BKEXToken)
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006404"
$ws.Range("E41").Value = "40KickTokenKICK"
# Row 42 - now BKEXToken (was CEJI)
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1050"
$ws.Range("E42").Value = "41BKEXTokenBKK"
# Row 43 - now CEJI (was KickToken)
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002246"
$ws.Range("E43").Value = "42CEJICEJI"

# Row 44 - LocalTraders
$ws.Range("D44").Value = "'0.008094"
# Row 45 - CoinLion
$ws.Range("D45").Value = "'0.00005475"
# Row 46 - Kangarootoken
$ws.Range("D46").Value = "'0.00000000750"
# Row 47 - CoinbaseStockToken
$ws.Range("D47").Value = "'0.4540"
# Row 48 - BOLO
$ws.Range("D48").Value = "'0.002565"
# Row 49 - CryptobidCoin
$ws.Range("D49").Value = "'0.00002101"
# Row 50 - SpecialPowerGold
$ws.Range("D50").Value = "'0.0002001"
